$wb = $excel.ActiveWorkbook

# --- Sheet: ELF-bldg-winter ---
$ws = $wb.Worksheets.Item("ELF-bldg-winter")
$ws.Range("B2").Value = 3.24256
$ws.Range("D2").Value = 2.79864
$ws.Range("B5").Value = 1.54842
$ws.Range("D5").Value = 1.0184
$ws.Range("B6").Value = 1.72105
$ws.Range("D6").Value = 1.33954
$ws.Range("B7").Value = 1.49589
$ws.Range("D7").Value = 1.0184

# --- Sheet: ELF-bldg-summer ---
$ws = $wb.Worksheets.Item("ELF-bldg-summer")
$ws.Range("B3").Value = 7.84772
$ws.Range("D3").Value = 4.66709
$ws.Range("B5").Value = 1.49329
$ws.Range("D5").Value = 1.67929
$ws.Range("B6").Value = 0.65358
$ws.Range("D6").Value = 1.52153
$ws.Range("B7").Value = 1.27475
$ws.Range("D7").Value = 1.67929

# --- Sheet: ELF-vehicles ---
$ws = $wb.Worksheets.Item("ELF-vehicles")
$ws.Range("B2").Value = 2.49052
$ws.Range("C2").Value = 2.33913
$ws.Range("B3").Value = 0.97658
$ws.Range("C3").Value = 0.9167999999999999
$ws.Range("B4").Value = 1.11083
$ws.Range("C4").Value = 0.96409
$ws.Range("B5").Value = 1.11083
$ws.Range("C5").Value = 0.96409
$ws.Range("B6").Value = 1.11083
$ws.Range("C6").Value = 0.96409
$ws.Range("B7").Value = 1.11083
$ws.Range("C7").Value = 0.96409

# --- Sheet: ELF-sectors ---
$ws = $wb.Worksheets.Item("ELF-sectors")
$ws.Range("B6").Value = 1.53238
$ws.Range("C6").Value = 1.0979

$wb.Save()
